$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1275699870682362
$ws.Range("C2").Value = 0.2949975552939468
$ws.Range("D2").Value = 0.150653331260784
$ws.Range("E2").Value = 0.3881408652290866
$ws.Range("F2").Value = 0.3748164570332126
$ws.Range("B3").Value = 0.0924332518450889
$ws.Range("C3").Value = 0.3674313087794084
$ws.Range("D3").Value = 0.3855150429668173
$ws.Range("E3").Value = 0.620898577037198
$ws.Range("F3").Value = 0.6277787736124477
$ws.Range("B4").Value = 0.6414105736191978
$ws.Range("C4").Value = 0.7989529440794614
$ws.Range("D4").Value = 4.695416189738523
$ws.Range("E4").Value = 2.166890903977061
$ws.Range("F4").Value = 2.116302514390309
$ws.Range("B5").Value = 0.2158879425535663
$ws.Range("C5").Value = 1.111983560722235
$ws.Range("D5").Value = 7.608038160243865
$ws.Range("E5").Value = 2.758267238728667
$ws.Range("F5").Value = 2.811606648419342
$ws.Range("B6").Value = 0.248349034037972
$ws.Range("C6").Value = 1.124424048567344
$ws.Range("D6").Value = 7.646948434845032
$ws.Range("E6").Value = 2.765311634309058
$ws.Range("F6").Value = 2.81603555487595
$ws.Range("B7").Value = 0.299637262814125
$ws.Range("C7").Value = 1.183108062718223
$ws.Range("D7").Value = 7.779390080635538
$ws.Range("E7").Value = 2.789155800710232
$ws.Range("F7").Value = 2.835336877543529
$ws.Range("B8").Value = 0.2652960715237629
$ws.Range("C8").Value = 1.272191104739367
$ws.Range("D8").Value = 7.927883031278751
$ws.Range("E8").Value = 2.815649664159011
$ws.Range("F8").Value = 2.866122845107168
$ws.Range("B9").Value = 0.1733904757308894
$ws.Range("C9").Value = 1.27878032370931
$ws.Range("D9").Value = 7.990910996653001
$ws.Range("E9").Value = 2.826819944151555
$ws.Range("F9").Value = 2.884909601627003
$ws.Range("B10").Value = 0.2219271402009887
$ws.Range("C10").Value = 1.341057205553475
$ws.Range("D10").Value = 8.037643935850472
$ws.Range("E10").Value = 2.835073885430585
$ws.Range("F10").Value = 2.889896368437108
$ws.Range("B11").Value = 0.1895377532140145
$ws.Range("C11").Value = 1.293575684359516
$ws.Range("D11").Value = 8.015790308611123
$ws.Range("E11").Value = 2.831217107289924
$ws.Range("F11").Value = 2.888353669326865

Write-Output "Updated forecast error values for Component Analysis filtering."
